$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New FF (FeedForward) data rows gathered 20-09-19, appended after the existing
# data (which currently ends at row 43).
$newRows = @(
    @(30000.0, 10.0, 3.0, 0.005, 200.0, 1.0,    0.30176211453744495, 2.0,   $false),
    @(30000.0, 10.0, 3.0, 0.005, 200.0, 5.0,    0.2472466960352423,  5.0,   $false),
    @(30000.0, 10.0, 3.0, 0.005, 200.0, 10.0,   0.5275330396475771,  8.0,   $false),
    @(30000.0, 10.0, 3.0, 0.005, 200.0, 20.0,   0.6227973568281938,  14.0,  $false),
    @(30000.0, 10.0, 3.0, 0.005, 200.0, 50.0,   0.8012114537444934,  33.0,  $false),
    @(30000.0, 10.0, 3.0, 0.005, 200.0, 100.0,  0.8700440528634361,  64.0,  $false),
    @(30000.0, 10.0, 3.0, 0.005, 200.0, 500.0,  0.9416299559471366,  316.0, $false),
    @(30000.0, 10.0, 3.0, 0.005, 200.0, 1000.0, 0.9564977973568282,  490.0, $false)
)

$startRow = 44
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
